# Update odds values in row 10 (Millonarios - Santa Fe match) on Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G10").Value = 2.15
$ws.Range("H10").Value = 2.8
$ws.Range("I10").Value = 4.2
$ws.Range("M10").Value = 1.14
$ws.Range("N10").Value = 5.5
$ws.Range("U10").Value = 2.38
$ws.Range("V10").Value = 1.53
$ws.Range("X10").Value = 8.5
$ws.Range("AC10").Value = 5.5
$ws.Range("AE10").Value = 21
$ws.Range("AF10").Value = 81
$ws.Range("AK10").Value = 41
$ws.Range("AN10").Value = 4
$ws.Range("AQ10").Value = 51
$ws.Range("AV10").Value = 81
$ws.Range("AX10").Value = 5.5
